$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555378532143"
$ws1.Range("B2").Value = "go_stims-16512555378162155.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555378362148.csv"
$ws1.Range("B4").Value = "go_stims-16512555378382144.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555378522124.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651255540253471"
$ws2.Range("B2").Value = "ZB-match_4-16512555381126966.csv"
$ws2.Range("B3").Value = "OB-16512555383436956.csv"
$ws2.Range("B4").Value = "OB-1651255539184693.csv"
$ws2.Range("B5").Value = "TB-16512555396584706.csv"
$ws2.Range("B6").Value = "ZB-match_3-16512555383156946.csv"
$ws2.Range("B7").Value = "TB-16512555399544728.csv"
$ws2.Range("B8").Value = "OB-16512555385876927.csv"
$ws2.Range("B9").Value = "TB-16512555402334714.csv"
$ws2.Range("B10").Value = "ZB-match_6-16512555381956995.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555402544723"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555403014758"
$ws4.Range("B2").Value = "MM_stims-16512555402694707.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555402574785.csv"
$ws4.Range("B4").Value = "MM_stims-16512555402854736.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555402704742.csv"
$ws4.Range("B6").Value = "MM_stims-16512555403004737.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555402864735.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1651255540380471"
$ws5.Range("B2").Value = "SAT_stims-16512555403334792.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555403484719.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555403644712.csv"
$ws5.Range("B5").Value = "SAT_stims-1651255540307475.csv"
